{"js": "// Insert two new \"caption text\" paragraphs (styles \"Captioned Figure\" and\n// \"Image Caption\") right after the BodyText paragraph that ends with\n// \"...combined population of LA and Chicago.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its distinctive trailing text.\nconst anchorText = \"combined population of LA and Chicago.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(anchorText) !== -1) {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph ending in '\" + anchorText + \"'\");\n}\n\n// Insert the \"Captioned Figure\" caption paragraph right after the anchor.\nconst captionedFigurePara = anchor.insertParagraph(\"caption text\", Word.InsertLocation.after);\ncaptionedFigurePara.style = \"Captioned Figure\";\n\n// Insert the \"Image Caption\" caption paragraph right after that one.\nconst imageCaptionPara = captionedFigurePara.insertParagraph(\"caption text\", Word.InsertLocation.after);\nimageCaptionPara.style = \"Image Caption\";\n\nawait context.sync();\n", "ps1": "# Insert two new \"caption text\" paragraphs (styles \"Captioned Figure\" and\n# \"Image Caption\") right after the BodyText paragraph that ends with\n# \"...combined population of LA and Chicago.\"\n\n$d = $word.ActiveDocument\n\n$anchorText = \"combined population of LA and Chicago.\"\n$anchor = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find anchor paragraph ending in '$anchorText'\"\n}\n\n# Insert the \"Captioned Figure\" caption paragraph right after the anchor.\n$anchor.Range.InsertParagraphAfter()\n$captionedFigurePara = $anchor.Next()\n$captionedFigurePara.Range.Text = \"caption text\"\n$captionedFigurePara.Style = \"Captioned Figure\"\n\n# Insert the \"Image Caption\" caption paragraph right after that one.\n$captionedFigurePara.Range.InsertParagraphAfter()\n$imageCaptionPara = $captionedFigurePara.Next()\n$imageCaptionPara.Range.Text = \"caption text\"\n$imageCaptionPara.Style = \"Image Caption\"\n"}
